# nano-metadata.xlsx regeneration
#  1. Rename the first sheet ("Export this as TSV" -> "Export as TSV")
#  2. Freeze the header row on that sheet
#  3. Add ErrorTitle/ErrorMessage text to the existing data validations

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Rename sheet ------------------------------------------------------
$ws.Name = "Export as TSV"

# --- 2. Freeze top row -----------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Data validation error messages -------------------------------------
$ws.Range("I2:I1048576").Validation.ErrorTitle = "Value must come from list"
$ws.Range("I2:I1048576").Validation.ErrorMessage = "Value must be one of: mass_spectrometry_imaging."

$ws.Range("J2:J1048576").Validation.ErrorTitle = "Value must come from list"
$ws.Range("J2:J1048576").Validation.ErrorMessage = "Value must be one of: NanoDESI / NanoPOTS."

$ws.Range("K2:K1048576").Validation.ErrorTitle = "Value must come from list"
$ws.Range("K2:K1048576").Validation.ErrorMessage = "Value must be one of: metabolites_and_lipids."

$ws.Range("L2:L1048576").Validation.ErrorTitle = "Not a boolean"
$ws.Range("L2:L1048576").Validation.ErrorMessage = 'The values in this column must be "TRUE" or "FALSE".'

$ws.Range("O2:O1048576").Validation.ErrorTitle = "Value must come from list"
$ws.Range("O2:O1048576").Validation.ErrorMessage = "Value must be one of: MALDI / MALDI-2 / DESI / SIMS / nESI."

$ws.Range("P2:P1048576").Validation.ErrorTitle = "Value must come from list"
$ws.Range("P2:P1048576").Validation.ErrorMessage = "Value must be one of: MALDI / MALDI-2 / DESI / SIMS / nESI."

$ws.Range("Q2:Q1048576").Validation.ErrorTitle = "Not a number"
$ws.Range("Q2:Q1048576").Validation.ErrorMessage = "The values in this column must be numbers."

$ws.Range("R2:R1048576").Validation.ErrorTitle = "Not a number"
$ws.Range("R2:R1048576").Validation.ErrorMessage = "The values in this column must be numbers."

$ws.Range("S2:S1048576").Validation.ErrorTitle = "Not a number"
$ws.Range("S2:S1048576").Validation.ErrorMessage = "The values in this column must be numbers."

$ws.Range("T2:T1048576").Validation.ErrorTitle = "Value must come from list"
$ws.Range("T2:T1048576").Validation.ErrorMessage = "Value must be one of: nm / um."

$ws.Range("U2:U1048576").Validation.ErrorTitle = "Not a number"
$ws.Range("U2:U1048576").Validation.ErrorMessage = "The values in this column must be numbers."

$ws.Range("V2:V1048576").Validation.ErrorTitle = "Value must come from list"
$ws.Range("V2:V1048576").Validation.ErrorMessage = "Value must be one of: nm / um."
